$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove hidden chart-tracking defined names (_xlchart.v1.*)
for ($i = $wb.Names.Count; $i -ge 1; $i--) {
    $wb.Names.Item($i).Delete()
}

# Update recalculated evaluation results (values changed due to re-run / correction)
$ws.Range("B8").Value = 29.83838384
$ws.Range("C8").Value = 171.0909091
$ws.Range("D8").Value = 141.2525253
$ws.Range("E8").Value = 2.303030303
$ws.Range("F8").Value = 6.919191919
$ws.Range("G8").Value = 4.743379525
$ws.Range("H8").Value = 1.371770238
$ws.Range("I8").Value = 2.071374274
$ws.Range("B9").Value = 35.50505051
$ws.Range("C9").Value = 210.2626263
$ws.Range("D9").Value = 174.7575758
$ws.Range("E9").Value = 2.292929293
$ws.Range("F9").Value = 7.313131313
$ws.Range("G9").Value = 4.929735512
$ws.Range("H9").Value = 1.460244136
$ws.Range("I9").Value = 2.37592687
$ws.Range("B10").Value = 41.32
$ws.Range("C10").Value = 252.12
$ws.Range("D10").Value = 210.8
$ws.Range("E10").Value = 2.27
$ws.Range("F10").Value = 7.69
$ws.Range("G10").Value = 5.103409906
$ws.Range("H10").Value = 1.532455846
$ws.Range("I10").Value = 2.559488589
$ws.Range("B11").Value = 46.80412371
$ws.Range("C11").Value = 295.9072165
$ws.Range("D11").Value = 249.1030928
$ws.Range("E11").Value = 2.340206186
$ws.Range("F11").Value = 8.24742268
$ws.Range("G11").Value = 5.321230226
$ws.Range("H11").Value = 1.647609853
$ws.Range("I11").Value = 3.02126677
$ws.Range("B12").Value = 52.79381443
$ws.Range("C12").Value = 344.6082474
$ws.Range("D12").Value = 291.814433
$ws.Range("E12").Value = 2.319587629
$ws.Range("F12").Value = 8.75257732
$ws.Range("G12").Value = 5.531400184
$ws.Range("H12").Value = 1.788056548
$ws.Range("I12").Value = 3.538067093
$ws.Range("K14").Value = 70.52525253
$ws.Range("L14").Value = 2249.141414
$ws.Range("M14").Value = 2178.616162
$ws.Range("N14").Value = 12.21212121
$ws.Range("O14").Value = 43.19191919
$ws.Range("P14").Value = 31.84380716
$ws.Range("Q14").Value = 8.344739394
$ws.Range("R14").Value = 79.04343952
$ws.Range("K15").Value = 111.2857143
$ws.Range("L15").Value = 5308.397959
$ws.Range("M15").Value = 5197.112245
$ws.Range("N15").Value = 14.64285714
$ws.Range("O15").Value = 64.90816327
$ws.Range("P15").Value = 48.43510155
$ws.Range("Q15").Value = 13.31927032
$ws.Range("R15").Value = 190.6862272
$ws.Range("K16").Value = 159.8125
$ws.Range("L16").Value = 9403.427083
$ws.Range("M16").Value = 9243.614583
$ws.Range("N16").Value = 16.20833333
$ws.Range("O16").Value = 81.32291667
$ws.Range("P16").Value = 60.08271752
$ws.Range("Q16").Value = 16.98016702
$ws.Range("R16").Value = 308.8618536
$ws.Range("K17").Value = 197.3814433
$ws.Range("L17").Value = 15126.3299
$ws.Range("M17").Value = 14928.94845
$ws.Range("N17").Value = 19.82474227
$ws.Range("O17").Value = 100.8659794
$ws.Range("P17").Value = 77.9865941
$ws.Range("Q17").Value = 20.33635647
$ws.Range("R17").Value = 447.9300224
$ws.Range("B18").Value = 11.51515152
$ws.Range("C18").Value = 307.4949495
$ws.Range("D18").Value = 295.979798
$ws.Range("E18").Value = 15.91919192
$ws.Range("F18").Value = 33.27272727
$ws.Range("G18").Value = 26.21387058
$ws.Range("H18").Value = 6.345208558
$ws.Range("I18").Value = 50.76347283
$ws.Range("K18").Value = 242.7979798
$ws.Range("L18").Value = 21825.68687
$ws.Range("M18").Value = 21582.88889
$ws.Range("N18").Value = 20.97979798
$ws.Range("O18").Value = 118.010101
$ws.Range("P18").Value = 91.19087632
$ws.Range("Q18").Value = 23.80941337
$ws.Range("R18").Value = 607.1622597
$ws.Range("B19").Value = 14.02020202
$ws.Range("C19").Value = 414.979798
$ws.Range("D19").Value = 400.959596
$ws.Range("E19").Value = 15.96969697
$ws.Range("F19").Value = 37.17171717
$ws.Range("G19").Value = 29.20919943
$ws.Range("H19").Value = 7.201167031
$ws.Range("I19").Value = 60.43683608
$ws.Range("K19").Value = 291.1717172
$ws.Range("L19").Value = 29916.54545
$ws.Range("M19").Value = 29625.37374
$ws.Range("N19").Value = 20.27272727
$ws.Range("O19").Value = 135.969697
$ws.Range("P19").Value = 103.5416339
$ws.Range("Q19").Value = 28.87701839
$ws.Range("R19").Value = 884.4372901
$ws.Range("B20").Value = 16.3
$ws.Range("C20").Value = 543.62
$ws.Range("D20").Value = 527.32
$ws.Range("E20").Value = 17.21
$ws.Range("F20").Value = 42.38
$ws.Range("G20").Value = 33.14454963
$ws.Range("H20").Value = 8.552715923
$ws.Range("I20").Value = 88.03060109
$ws.Range("K20").Value = 339.36
$ws.Range("L20").Value = 38953.08
$ws.Range("M20").Value = 38613.72
$ws.Range("N20").Value = 17.77
$ws.Range("O20").Value = 151.8
$ws.Range("P20").Value = 115.8220744
$ws.Range("Q20").Value = 33.60365438
$ws.Range("R20").Value = 1199.486479
$ws.Range("B21").Value = 18.18556701
$ws.Range("C21").Value = 692.6597938
$ws.Range("D21").Value = 674.4742268
$ws.Range("E21").Value = 19.18556701
$ws.Range("F21").Value = 47.92783505
$ws.Range("G21").Value = 37.70811621
$ws.Range("H21").Value = 9.398153753
$ws.Range("I21").Value = 100.1391471
$ws.Range("K21").Value = 374.1237113
$ws.Range("L21").Value = 49206.16495
$ws.Range("M21").Value = 48832.04124
$ws.Range("N21").Value = 21.95876289
$ws.Range("O21").Value = 172.5670103
$ws.Range("P21").Value = 132.5040794
$ws.Range("Q21").Value = 37.05771867
$ws.Range("R21").Value = 1424.524485
$ws.Range("B22").Value = 20.30927835
$ws.Range("C22").Value = 846.8556701
$ws.Range("D22").Value = 826.5463918
$ws.Range("E22").Value = 20.05154639
$ws.Range("F22").Value = 52.88659794
$ws.Range("G22").Value = 41.39290198
$ws.Range("H22").Value = 10.51740847
$ws.Range("I22").Value = 123.9160647
$ws.Range("K22").Value = 410.742268
$ws.Range("L22").Value = 60834.75258
$ws.Range("M22").Value = 60424.01031
$ws.Range("N22").Value = 23.29896907
$ws.Range("O22").Value = 193.1958763
$ws.Range("P22").Value = 149.1110428
$ws.Range("Q22").Value = 41.5851288
$ws.Range("R22").Value = 1801.413258
$ws.Range("C28").Value = 351.6868687
$ws.Range("D28").Value = 350.6868687
$ws.Range("E28").Value = 350.6868687
$ws.Range("F28").Value = 350.6868687
$ws.Range("G28").Value = 350.6868687
$ws.Range("C29").Value = 368.9393939
$ws.Range("D29").Value = 367.9393939
$ws.Range("E29").Value = 367.9393939
$ws.Range("F29").Value = 367.9393939
$ws.Range("G29").Value = 367.9393939
$ws.Range("C30").Value = 401.67
$ws.Range("D30").Value = 400.67
$ws.Range("E30").Value = 400.67
$ws.Range("F30").Value = 400.67
$ws.Range("G30").Value = 400.67
$ws.Range("C31").Value = 425.3917526
$ws.Range("D31").Value = 424.3917526
$ws.Range("E31").Value = 424.3917526
$ws.Range("F31").Value = 424.3917526
$ws.Range("G31").Value = 424.3917526
$ws.Range("C32").Value = 437.6391753
$ws.Range("D32").Value = 436.6391753
$ws.Range("E32").Value = 436.6391753
$ws.Range("F32").Value = 436.6391753
$ws.Range("G32").Value = 436.6391753
$ws.Range("C38").Value = 65.64646465
$ws.Range("D38").Value = 64.64646465
$ws.Range("E38").Value = 64.64646465
$ws.Range("F38").Value = 64.64646465
$ws.Range("G38").Value = 64.64646465
$ws.Range("C39").Value = 76.85858586
$ws.Range("D39").Value = 75.85858586
$ws.Range("E39").Value = 75.85858586
$ws.Range("F39").Value = 75.85858586
$ws.Range("G39").Value = 75.85858586
$ws.Range("C40").Value = 88.94
$ws.Range("D40").Value = 87.94
$ws.Range("E40").Value = 87.94
$ws.Range("F40").Value = 87.94
$ws.Range("G40").Value = 87.94
$ws.Range("C41").Value = 99.17525773
$ws.Range("D41").Value = 98.17525773
$ws.Range("E41").Value = 98.17525773
$ws.Range("F41").Value = 98.17525773
$ws.Range("G41").Value = 98.17525773
$ws.Range("C42").Value = 109.0103093
$ws.Range("D42").Value = 108.0103093
$ws.Range("E42").Value = 108.0103093
$ws.Range("F42").Value = 108.0103093
$ws.Range("G42").Value = 108.0103093

# Restore selection/view state to match the edited workbook
$ws.Range("K34:R42").Select()
